# Sprint Backlog Template - reassign the "Assignee" column for two user
# stories:
#   D3  (1.1 "...create an account.")      Ahmed Ashraf -> Bishoy
#   D11 (2.4 "...delete a task.")          Bishoy       -> Shamel
#
# Also normalize the font used on D4/D6 (which previously carried the
# non-standard "Docs-Calibri" font name left over from the source data)
# back to the workbook's default Calibri font, matching the rest of the
# Assignee column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "Bishoy"
$ws.Range("D11").Value = "Shamel"

$ws.Range("D4").Font.Name = "Calibri"
$ws.Range("D6").Font.Name = "Calibri"
